$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'The applicant''s projects demonstrate strong skills in ReactJS, Django, Python, and JavaScript, all of which are highly relevant to the job description. The applicant has experience in developing web applications and backend systems, displaying proficiency in MongoDB, ExpressJS, and NodeJS. Additionally, the applicant has worked on projects involving computer vision and image processing, showcasing a diverse skill set. However, the applicant lacks direct experience with MongoDB, which is a required skill for the job. While the applicant''s projects demonstrate strong technical skills, the lack of direct MongoDB experience could be a potential area for further development.'
$ws.Range("E2").Value = 'The applicant has demonstrated a willingness to learn, grow, and adapt through proactive participation in workshops and online courses to enhance skills. The applicant also exhibits a strong commitment to self-improvement and learning from setbacks, indicating good problem-solving skills and adaptability. The applicant''s expressed interest in Japanese work culture and values, as well as their desire to work in Japan, align well with the company''s cultural expectations. The applicant''s detailed expectations from the company and clear career plan reflect strong communication and time management skills. Overall, the applicant''s responses indicate a positive attitude, adaptability, and a strong alignment with the company''s soft skill requirements.'
$ws.Range("C3").Value = 'The applicant has a high level of skill in ReactJS, NodeJS, and WebRTC, which are all relevant to the job requirements. Although the applicant does not have direct experience with MongoDB, the candidate''s experience with Firebase, Docker, and ExpressJS demonstrates a strong foundation in database management and backend development, which can be transferrable to working with MongoDB. The applicant''s project experience in developing a time-saving food delivery platform and a web app for real-time gaming and video calling shows the ability to work on web and mobile applications, which aligns with the responsibilities of the SDE Intern position. The use of technologies like Flutter, Dart, and ReactJS in the ''SHAMIYANA APP'' project indicates the ability to work on front-end development. Overall, the applicant''s diverse skill set and experience make them well-suited for the role, earning a high score.'
$ws.Range("E3").Value = 'Based on the applicant''s answers, it is evident that the applicant possesses strong communication skills as reflected in the clear articulation of career plans and expectations from the company. The applicant also expresses a willingness to adapt to new cultures and an interest in learning new skills, illustrating strong adaptability. Additionally, the applicant''s acknowledgment of strengths in teamwork and open communication, along with the ability to address challenges collaboratively, demonstrates good problem-solving and teamwork skills. Overall, the applicant''s responses reflect a positive attitude, strong teamwork skills, and adaptability, warranting a score of 4.'
$ws.Range("C4").Value = 'The applicant has worked on several projects showcasing strong technical skills, including computer vision, image processing, and backend development using Django Rest Framework and ReactJS. While the specific skills required by the company, such as MongoDB and NodeJS, are not directly demonstrated in the projects, the applicant''s broad skill set and experience with related technologies make them well-equipped to learn and adapt to the company''s tech stack. The applicant''s experience in predicting human activities using data analysis and machine learning aligns with the company''s focus on web/mobile applications and product enhancement suggestions, demonstrating their potential to contribute effectively to the role.'
$ws.Range("E4").Value = 'The applicant has exhibited strong communication skills, expressing eagerness to immerse themselves in Japan''s AI industry and adapt to a new culture. Their interest in contributing to cutting-edge AI/ML technologies and desire to learn new skills align with the company''s emphasis on problem solving and adaptability. The applicant''s ability to recognize the importance of collaboration and open communication, as evidenced by their experience with project setbacks, indicates their potential to excel in a team-oriented environment. Given their positive attitude towards adapting to different cultures and their clear career goals, the applicant demonstrates a high level of enthusiasm and adaptability, earning a score of 4 out of 5 for their willingness to work in Japan.'
$ws.Range("C5").Value = 'The applicant has good experience working on web development projects such as the Video Conferencing Project, TEDxIITGuwahati Website, and the Alcheringa Pass Portal, which involved JavaScript, HTML, and CSS. These projects demonstrate the applicant''s ability to work with front-end technologies and API integration. Although the applicant''s projects do not directly align with the MERN stack, the skills obtained from these projects are transferrable and can be utilized in the internship. The applicant has experience in API integration, JavaScript, HTML, and CSS, which are crucial for web development. Overall, the applicant''s projects showcase a strong foundation in web development and the ability to work with relevant technologies, indicating a good potential fit for the MERN Stack Developer Internship. However, the applicant could benefit from further experience with MongoDB, ExpressJS, and NodeJS to fully align with the job requirements.'
$ws.Range("E5").Value = 'Based on the applicant''s responses, it is evident that they possess excellent communication skills and demonstrate a willingness to adapt to new cultures. The applicant''s interest in being part of Japan''s work culture and their openness to working in various countries showcase a strong level of adaptability. Additionally, the applicant has emphasized the importance of teamwork and collaboration, reflecting positive attributes that align with the company''s soft skill requirements. The applicant''s proactive approach to addressing weaknesses and their ability to learn from challenges further illustrate their problem-solving skills. Overall, the applicant exhibits strong attributes that align with the company''s soft skill requirements and shows a high level of willingness to work in Japan.'
